$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared strings: "MuSCs" cluster label becomes "Resolving-Mac",
# and the former "Inflammatory-Mac" cluster label becomes "MuSCs".
# Update every cell that referenced those labels (sending/target cluster columns).
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("D13").Value = "Resolving-Mac"

$ws.Range("D4").Value = "MuSCs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("D12").Value = "MuSCs"

# Updated NATMI TPM-derived statistics (columns F-T).
$ws.Range("G2").Value = 0.06762866666666667
$ws.Range("H2").Value = 0.202886
$ws.Range("I2").Value = 0.0134153952845566
$ws.Range("J2").Value = 0.0134153952845566
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2158443333333333
$ws.Range("N2").Value = 0.647533
$ws.Range("O2").Value = 0.0304269596471442
$ws.Range("P2").Value = 0.0383225350174662
$ws.Range("Q2").Value = 0.01459726447088889
$ws.Range("R2").Value = 0.131375380238
$ws.Range("S2").Value = 0.0004081896909736922
$ws.Range("T2").Value = 0.0005141119555655712
$ws.Range("G3").Value = 0.06762866666666667
$ws.Range("H3").Value = 0.202886
$ws.Range("I3").Value = 0.0134153952845566
$ws.Range("J3").Value = 0.0134153952845566
$ws.Range("O3").Value = 0.3272075825100088
$ws.Range("P3").Value = 0.4121155772426048
$ws.Range("Q3").Value = 0.1569770911773333
$ws.Range("R3").Value = 1.412793820596
$ws.Range("S3").Value = 0.004389619059475936
$ws.Range("T3").Value = 0.005528693371632761
$ws.Range("G4").Value = 0.06762866666666667
$ws.Range("H4").Value = 0.202886
$ws.Range("I4").Value = 0.0134153952845566
$ws.Range("J4").Value = 0.0134153952845566
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.3846295
$ws.Range("N4").Value = 8.769259
$ws.Range("O4").Value = 0.6180887068188557
$ws.Range("P4").Value = 0.5189854958816471
$ws.Range("Q4").Value = 0.2965266469123333
$ws.Range("R4").Value = 1.779159881474
$ws.Range("S4").Value = 0.008291904322895364
$ws.Range("T4").Value = 0.006962395574203917
$ws.Range("G5").Value = 0.06762866666666667
$ws.Range("H5").Value = 0.202886
$ws.Range("I5").Value = 0.0134153952845566
$ws.Range("J5").Value = 0.0134153952845566
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1722156666666667
$ws.Range("N5").Value = 0.516647
$ws.Range("O5").Value = 0.02427675102399122
$ws.Range("P5").Value = 0.03057639185828191
$ws.Range("Q5").Value = 0.01164671591577778
$ws.Range("R5").Value = 0.104820443242
$ws.Range("S5").Value = 0.0003256822112116064
$ws.Range("T5").Value = 0.0004101943831543498
$ws.Range("I6").Value = 0.9827953701592058
$ws.Range("J6").Value = 0.9827953701592059
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2158443333333333
$ws.Range("N6").Value = 0.647533
$ws.Range("O6").Value = 0.0304269596471442
$ws.Range("P6").Value = 0.0383225350174662
$ws.Range("Q6").Value = 1.069377654156333
$ws.Range("R6").Value = 9.624398887406999
$ws.Range("S6").Value = 0.0299034750692343
$ws.Range("T6").Value = 0.03766320998792982
$ws.Range("I7").Value = 0.9827953701592058
$ws.Range("J7").Value = 0.9827953701592059
$ws.Range("O7").Value = 0.3272075825100088
$ws.Range("P7").Value = 0.4121155772426048
$ws.Range("S7").Value = 0.321578097171823
$ws.Range("T7").Value = 0.4050252812845206
$ws.Range("I8").Value = 0.9827953701592058
$ws.Range("J8").Value = 0.9827953701592059
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.3846295
$ws.Range("N8").Value = 8.769259
$ws.Range("O8").Value = 0.6180887068188557
$ws.Range("P8").Value = 0.5189854958816471
$ws.Range("Q8").Value = 21.7231777023935
$ws.Range("R8").Value = 130.339066214361
$ws.Range("S8").Value = 0.6074547194092622
$ws.Range("T8").Value = 0.5100565425322624
$ws.Range("I9").Value = 0.9827953701592058
$ws.Range("J9").Value = 0.9827953701592059
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1722156666666667
$ws.Range("N9").Value = 0.516647
$ws.Range("O9").Value = 0.02427675102399122
$ws.Range("P9").Value = 0.03057639185828191
$ws.Range("Q9").Value = 0.8532240934236666
$ws.Range("R9").Value = 7.679016840812999
$ws.Range("S9").Value = 0.02385907850888633
$ws.Range("T9").Value = 0.0300503363544931
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.019102
$ws.Range("H10").Value = 0.057306
$ws.Range("I10").Value = 0.003789234556237495
$ws.Range("J10").Value = 0.003789234556237496
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2158443333333333
$ws.Range("N10").Value = 0.647533
$ws.Range("O10").Value = 0.0304269596471442
$ws.Range("P10").Value = 0.0383225350174662
$ws.Range("Q10").Value = 0.004123058455333333
$ws.Range("R10").Value = 0.03710752609800001
$ws.Range("S10").Value = 0.0001152948869362026
$ws.Range("T10").Value = 0.0001452130739708044
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.019102
$ws.Range("H11").Value = 0.057306
$ws.Range("I11").Value = 0.003789234556237495
$ws.Range("J11").Value = 0.003789234556237496
$ws.Range("O11").Value = 0.3272075825100088
$ws.Range("P11").Value = 0.4121155772426048
$ws.Range("Q11").Value = 0.044338836524
$ws.Range("R11").Value = 0.399049528716
$ws.Range("S11").Value = 0.001239866278709857
$ws.Range("T11").Value = 0.001561602586451441
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.019102
$ws.Range("H12").Value = 0.057306
$ws.Range("I12").Value = 0.003789234556237495
$ws.Range("J12").Value = 0.003789234556237496
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.3846295
$ws.Range("N12").Value = 8.769259
$ws.Range("O12").Value = 0.6180887068188557
$ws.Range("P12").Value = 0.5189854958816471
$ws.Range("Q12").Value = 0.083755192709
$ws.Range("R12").Value = 0.502531156254
$ws.Range("S12").Value = 0.002342083086698154
$ws.Range("T12").Value = 0.00196655777518079
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.019102
$ws.Range("H13").Value = 0.057306
$ws.Range("I13").Value = 0.003789234556237495
$ws.Range("J13").Value = 0.003789234556237496
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1722156666666667
$ws.Range("N13").Value = 0.516647
$ws.Range("O13").Value = 0.02427675102399122
$ws.Range("P13").Value = 0.03057639185828191
$ws.Range("Q13").Value = 0.003289663664666667
$ws.Range("R13").Value = 0.029606972982
$ws.Range("S13").Value = 0.00009199030389328153
$ws.Range("T13").Value = 0.0001158611206344606
